$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 6766.643
$ws.Range("I86").Value = 7338.381
$ws.Range("J86").Value = 5051.4287
$ws.Range("K86").Value = 7338.381
$ws.Range("L86").Value = 5051.4287
$ws.Range("M86").Value = -6215.381
$ws.Range("N86").Value = -7297.4287

$ws.Range("H88").Value = 3000
$ws.Range("I88").Value = 1000
$ws.Range("J88").Value = 5000
$ws.Range("K88").Value = 1000
$ws.Range("L88").Value = 5000
$ws.Range("M88").Value = -594
$ws.Range("N88").Value = -5812

$ws.Range("H89").Value = 6766.643
$ws.Range("I89").Value = 7338.381
$ws.Range("J89").Value = 5051.4287
$ws.Range("K89").Value = 36691.905
$ws.Range("L89").Value = 25257.1435
$ws.Range("M89").Value = -31075.905
$ws.Range("N89").Value = -36489.14350000001

$ws.Range("H91").Value = 3000
$ws.Range("I91").Value = 1000
$ws.Range("J91").Value = 5000
$ws.Range("K91").Value = 1000
$ws.Range("L91").Value = 5000
$ws.Range("M91").Value = 404
$ws.Range("N91").Value = -7808

$ws.Range("H96").Value = 297.90475
$ws.Range("I96").Value = 237.1579
$ws.Range("K96").Value = 711.4737
$ws.Range("M96").Value = 661.5263

$ws.Range("H132").Value = 2131.6462
$ws.Range("I132").Value = 1486.9857
$ws.Range("J132").Value = 5892.1665
$ws.Range("K132").Value = 4460.9571
$ws.Range("L132").Value = 17676.4995
$ws.Range("M132").Value = -1930.9571
$ws.Range("N132").Value = -22736.4995

$ws.Range("H138").Value = 3115.9587
$ws.Range("I138").Value = 1807.0333
$ws.Range("J138").Value = 3702.0447
$ws.Range("K138").Value = 5421.0999
$ws.Range("L138").Value = 11106.1341
$ws.Range("M138").Value = -281.0999000000002
$ws.Range("N138").Value = -21386.1341

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 21518.75
$ws.Range("I74").Value = 1773.7174
$ws.Range("J74").Value = 172897.33
$ws.Range("K74").Value = 1773.7174
$ws.Range("L74").Value = 172897.33
$ws.Range("M74").Value = -899.7174
$ws.Range("N74").Value = -174645.33

$ws.Range("H77").Value = 21518.75
$ws.Range("I77").Value = 1773.7174
$ws.Range("J77").Value = 172897.33
$ws.Range("K77").Value = 8868.587
$ws.Range("L77").Value = 864486.6499999999
$ws.Range("M77").Value = -4500.587
$ws.Range("N77").Value = -873222.6499999999

$ws.Range("H122").Value = 12814
$ws.Range("I122").Value = 14884.4
$ws.Range("J122").Value = 5050
$ws.Range("K122").Value = 44653.2
$ws.Range("L122").Value = 15150
$ws.Range("M122").Value = -42203.2
$ws.Range("N122").Value = -20050

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2692.9092
$ws.Range("I134").Value = 2216
$ws.Range("K134").Value = 6648
$ws.Range("M134").Value = -4113

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2104500.8
$ws.Range("I16").Value = 3272229.5
$ws.Range("J16").Value = 2589.2
$ws.Range("K16").Value = 3272229.5
$ws.Range("L16").Value = 2589.2
$ws.Range("M16").Value = -3271942.5
$ws.Range("N16").Value = -3163.2

$ws.Range("H31").Value = 4476.109
$ws.Range("I31").Value = 2424.923
$ws.Range("J31").Value = 5111
$ws.Range("K31").Value = 2424.923
$ws.Range("L31").Value = 5111
$ws.Range("M31").Value = -2129.923
$ws.Range("N31").Value = -5701

$ws.Range("H34").Value = 4476.109
$ws.Range("I34").Value = 2424.923
$ws.Range("J34").Value = 5111
$ws.Range("K34").Value = 2424.923
$ws.Range("L34").Value = 5111
$ws.Range("M34").Value = -2222.923
$ws.Range("N34").Value = -5515

$ws.Range("H113").Value = 2104500.8
$ws.Range("I113").Value = 3272229.5
$ws.Range("J113").Value = 2589.2
$ws.Range("K113").Value = 3272229.5
$ws.Range("L113").Value = 2589.2
$ws.Range("M113").Value = -3270059.5
$ws.Range("N113").Value = -6929.2

$ws.Range("H134").Value = 3051.353
$ws.Range("I134").Value = 3051.353
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 9154.059000000001
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -6619.059000000001
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 14272.728
$ws.Range("I9").Value = 3000
$ws.Range("J9").Value = 15400
$ws.Range("K9").Value = 9000
$ws.Range("L9").Value = 46200
$ws.Range("M9").Value = -8776
$ws.Range("N9").Value = -46648

$ws.Range("H16").Value = 4434.0586
$ws.Range("I16").Value = 2198.5
$ws.Range("J16").Value = 5121.923
$ws.Range("K16").Value = 6595.5
$ws.Range("L16").Value = 15365.769
$ws.Range("M16").Value = -6422.5
$ws.Range("N16").Value = -15711.769

$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()

$ws.Range("H22").Value = 986.1429000000001
$ws.Range("I22").Value = 545
$ws.Range("J22").Value = 1089.9412
$ws.Range("K22").Value = 1635
$ws.Range("L22").Value = 3269.8236
$ws.Range("M22").Value = -1466
$ws.Range("N22").Value = -3607.8236

$ws.Range("H27").Value = 986.1429000000001
$ws.Range("I27").Value = 545
$ws.Range("J27").Value = 1089.9412
$ws.Range("K27").Value = 1635
$ws.Range("L27").Value = 3269.8236
$ws.Range("M27").Value = -1533
$ws.Range("N27").Value = -3473.8236

$ws.Range("H32").Value = 2000
$ws.Range("I32").Value = 2000
$ws.Range("K32").Value = 6000
$ws.Range("M32").Value = -5717

$ws.Range("H64").Value = 9289.846
$ws.Range("J64").Value = 9972.333000000001
$ws.Range("L64").Value = 29916.999
$ws.Range("N64").Value = -30456.999

$ws.Range("H67").Value = 9289.846
$ws.Range("J67").Value = 9972.333000000001
$ws.Range("L67").Value = 29916.999
$ws.Range("N67").Value = -31788.999

$ws.Range("H70").Value = 6713.385
$ws.Range("J70").Value = 8388.888999999999
$ws.Range("L70").Value = 25166.667
$ws.Range("N70").Value = -25796.667

$ws.Range("H73").Value = 6713.385
$ws.Range("J73").Value = 8388.888999999999
$ws.Range("L73").Value = 25166.667
$ws.Range("N73").Value = -27350.667

$ws.Range("H76").Value = 7166.6665
$ws.Range("I76").Value = 10233.333
$ws.Range("J76").Value = 6400
$ws.Range("K76").Value = 30699.999
$ws.Range("L76").Value = 19200
$ws.Range("M76").Value = -30316.999
$ws.Range("N76").Value = -19966

$ws.Range("H79").Value = 7166.6665
$ws.Range("I79").Value = 10233.333
$ws.Range("J79").Value = 6400
$ws.Range("K79").Value = 30699.999
$ws.Range("L79").Value = 19200
$ws.Range("M79").Value = -29373.999
$ws.Range("N79").Value = -21852

$ws.Range("H103").Value = 893.2222
$ws.Range("I103").Value = 350
$ws.Range("J103").Value = 1001.86664
$ws.Range("K103").Value = 1050
$ws.Range("L103").Value = 3005.59992
$ws.Range("M103").Value = -171
$ws.Range("N103").Value = -4763.59992

$ws.Range("H109").Value = 2451.9443
$ws.Range("I109").Value = 756.1111
$ws.Range("K109").Value = 2268.3333
$ws.Range("M109").Value = -1228.3333

$ws.Range("H112").Value = 46353824
$ws.Range("I112").Value = 1500
$ws.Range("J112").Value = 49080430
$ws.Range("K112").Value = 4500
$ws.Range("L112").Value = 147241290
$ws.Range("M112").Value = -3392
$ws.Range("N112").Value = -147243506

$ws.Range("H122").Value = 596.5
$ws.Range("I122").Value = 306.48148
$ws.Range("J122").Value = 1466.5555
$ws.Range("K122").Value = 2758.33332
$ws.Range("L122").Value = 13198.9995
$ws.Range("M122").Value = -308.3333199999997
$ws.Range("N122").Value = -18098.9995

$ws.Range("H125").Value = 1783.1818
$ws.Range("I125").Value = 410
$ws.Range("K125").Value = 1230
$ws.Range("M125").Value = 3690

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 53548.81
$ws.Range("I70").Value = 58669.26
$ws.Range("J70").Value = 4904.5
$ws.Range("K70").Value = 58669.26
$ws.Range("L70").Value = 4904.5
$ws.Range("M70").Value = -58399.26
$ws.Range("N70").Value = -5444.5

$ws.Range("H73").Value = 53548.81
$ws.Range("I73").Value = 58669.26
$ws.Range("J73").Value = 4904.5
$ws.Range("K73").Value = 58669.26
$ws.Range("L73").Value = 4904.5
$ws.Range("M73").Value = -57733.26
$ws.Range("N73").Value = -6776.5

$ws.Range("H122").Value = 2646.28
$ws.Range("I122").Value = 2643.8
$ws.Range("J122").Value = 2650
$ws.Range("K122").Value = 7931.400000000001
$ws.Range("L122").Value = 7950
$ws.Range("M122").Value = -5481.400000000001
$ws.Range("N122").Value = -12850

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1645.6666
$ws.Range("I82").Value = 1274.1177
$ws.Range("J82").Value = 2277.3
$ws.Range("K82").Value = 1274.1177
$ws.Range("L82").Value = 2277.3
$ws.Range("M82").Value = -913.1177
$ws.Range("N82").Value = -2999.3

$ws.Range("H85").Value = 1645.6666
$ws.Range("I85").Value = 1274.1177
$ws.Range("J85").Value = 2277.3
$ws.Range("K85").Value = 1274.1177
$ws.Range("L85").Value = 2277.3
$ws.Range("M85").Value = -26.11770000000001
$ws.Range("N85").Value = -4773.3

$ws.Range("H107").Value = 2000
$ws.Range("I107").Value = 2000
$ws.Range("K107").Value = 2000
$ws.Range("M107").Value = -80

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H135").Value = 53832.273
$ws.Range("J135").Value = 53832.273
$ws.Range("L135").Value = 53832.273
$ws.Range("N135").Value = -63972.273

$ws.Range("H136").Value = 1048.1562
$ws.Range("I136").Value = 1039.9434
$ws.Range("J136").Value = 1087.7273
$ws.Range("K136").Value = 3119.8302
$ws.Range("L136").Value = 3263.1819
$ws.Range("M136").Value = -569.8302000000003
$ws.Range("N136").Value = -8363.1819

$ws.Range("H137").Value = 57571.332
$ws.Range("J137").Value = 57571.332
$ws.Range("L137").Value = 57571.332
$ws.Range("N137").Value = -67771.33199999999
